$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 989.34485
$ws.Range("J17").Value = 1091.2916
$ws.Range("L17").Value = 3273.8748
$ws.Range("N17").Value = -3609.8748

$ws.Range("H112").Value = 3272
$ws.Range("J112").Value = 3540.9092
$ws.Range("L112").Value = 10622.7276
$ws.Range("N112").Value = -12838.7276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 21400
$ws.Range("J24").Value = 21400
$ws.Range("L24").Value = 21400
$ws.Range("N24").Value = -22148

$ws.Range("H32").Value = 3044.44
$ws.Range("I32").Value = 3044.44
$ws.Range("K32").Value = 3044.44
$ws.Range("M32").Value = -2757.44

$ws.Range("H95").Value = 14177.462
$ws.Range("J95").Value = 14177.462
$ws.Range("L95").Value = 14177.462
$ws.Range("N95").Value = -19669.462

$ws.Range("H100").Value = 21400
$ws.Range("J100").Value = 21400
$ws.Range("L100").Value = 21400
$ws.Range("N100").Value = -23564

$ws.Range("H110").Value = 1352.4117
$ws.Range("I110").Value = 882.5833
$ws.Range("J110").Value = 2480
$ws.Range("K110").Value = 882.5833
$ws.Range("L110").Value = 2480
$ws.Range("M110").Value = 1162.4167
$ws.Range("N110").Value = -6570

$ws.Range("H114").Value = 15731.333
$ws.Range("J114").Value = 15731.333
$ws.Range("L114").Value = 15731.333
$ws.Range("N114").Value = -24409.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 13419.2
$ws.Range("I96").Value = 6952
$ws.Range("J96").Value = 23120
$ws.Range("K96").Value = 6952
$ws.Range("L96").Value = 23120
$ws.Range("M96").Value = -4206
$ws.Range("N96").Value = -28612

$ws.Range("H105").Value = 2289.2307
$ws.Range("I105").Value = 2128
$ws.Range("J105").Value = 2826.6667
$ws.Range("K105").Value = 2128
$ws.Range("L105").Value = 2826.6667
$ws.Range("M105").Value = -381
$ws.Range("N105").Value = -6320.6667

$ws.Range("H109").Value = 32000
$ws.Range("J109").Value = 32000
$ws.Range("L109").Value = 32000
$ws.Range("N109").Value = -34774

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2462.2446
$ws.Range("I31").Value = 1314.0769
$ws.Range("J31").Value = 3883.7856
$ws.Range("K31").Value = 1314.0769
$ws.Range("L31").Value = 3883.7856
$ws.Range("M31").Value = -1019.0769
$ws.Range("N31").Value = -4473.7856

$ws.Range("H34").Value = 2462.2446
$ws.Range("I34").Value = 1314.0769
$ws.Range("J34").Value = 3883.7856
$ws.Range("K34").Value = 1314.0769
$ws.Range("L34").Value = 3883.7856
$ws.Range("M34").Value = -1112.0769
$ws.Range("N34").Value = -4287.7856

$ws.Range("H58").Value = 3109.6978
$ws.Range("I58").Value = 1480.2
$ws.Range("J58").Value = 5372.8887
$ws.Range("K58").Value = 1480.2
$ws.Range("L58").Value = 5372.8887
$ws.Range("M58").Value = -1277.2
$ws.Range("N58").Value = -5778.8887

$ws.Range("H107").Value = 2024.1818
$ws.Range("I107").Value = 710
$ws.Range("J107").Value = 3601.2
$ws.Range("K107").Value = 710
$ws.Range("L107").Value = 3601.2
$ws.Range("M107").Value = 1210
$ws.Range("N107").Value = -7441.2

$ws.Range("H132").Value = 1721.6471
$ws.Range("I132").Value = 1165.0312
$ws.Range("J132").Value = 2659.1052
$ws.Range("K132").Value = 3495.0936
$ws.Range("L132").Value = 7977.3156
$ws.Range("M132").Value = -965.0935999999997
$ws.Range("N132").Value = -13037.3156

$ws.Range("H136").Value = 3109.6978
$ws.Range("I136").Value = 1480.2
$ws.Range("J136").Value = 5372.8887
$ws.Range("K136").Value = 4440.6
$ws.Range("L136").Value = 16118.6661
$ws.Range("M136").Value = -1890.6
$ws.Range("N136").Value = -21218.6661

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 19823
$ws.Range("J103").Value = 19823
$ws.Range("L103").Value = 19823
$ws.Range("N103").Value = -22167

$ws.Range("H132").Value = 2440.8845
$ws.Range("I132").Value = 1525.84
$ws.Range("J132").Value = 3288.1482
$ws.Range("K132").Value = 4577.52
$ws.Range("L132").Value = 9864.444600000001
$ws.Range("M132").Value = -2047.52
$ws.Range("N132").Value = -14924.4446

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 22558.654
$ws.Range("I132").Value = 26425.578
$ws.Range("J132").Value = 9173.154
$ws.Range("K132").Value = 79276.734
$ws.Range("L132").Value = 27519.462
$ws.Range("M132").Value = -76746.734
$ws.Range("N132").Value = -32579.462

$ws.Range("H136").Value = 1899.5
$ws.Range("I136").Value = 1776.6666
$ws.Range("J136").Value = 3005
$ws.Range("K136").Value = 5329.9998
$ws.Range("L136").Value = 9015
$ws.Range("M136").Value = -2779.9998
$ws.Range("N136").Value = -14115

$ws.Range("H139").Value = 30000
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws.Range("H140").Value = 250031500
$ws.Range("J140").Value = 250031500
$ws.Range("L140").Value = 250031500
$ws.Range("N140").Value = -250041860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 9003.727999999999
$ws.Range("J69").Value = 9003.727999999999
$ws.Range("L69").Value = 9003.727999999999
$ws.Range("N69").Value = -10501.728

$ws.Range("H72").Value = 9003.727999999999
$ws.Range("J72").Value = 9003.727999999999
$ws.Range("L72").Value = 27011.184
$ws.Range("N72").Value = -34499.18399999999

$ws.Range("H94").Value = 14818.571
$ws.Range("J94").Value = 14818.571
$ws.Range("L94").Value = 14818.571
$ws.Range("N94").Value = -16620.571

$ws.Range("H122").Value = 2718.6316
$ws.Range("I122").Value = 2321.25
$ws.Range("J122").Value = 3399.8572
$ws.Range("K122").Value = 6963.75
$ws.Range("L122").Value = 10199.5716
$ws.Range("M122").Value = -4513.75
$ws.Range("N122").Value = -15099.5716

$ws.Range("H132").Value = 2037.4634
$ws.Range("I132").Value = 1703.3182
$ws.Range("J132").Value = 2424.3684
$ws.Range("K132").Value = 5109.9546
$ws.Range("L132").Value = 7273.1052
$ws.Range("M132").Value = -2579.9546
$ws.Range("N132").Value = -12333.1052

$ws.Range("H136").Value = 1954.415
$ws.Range("I136").Value = 1966.9231
$ws.Range("J136").Value = 1919.5714
$ws.Range("K136").Value = 5900.7693
$ws.Range("L136").Value = 5758.7142
$ws.Range("M136").Value = -3350.7693
$ws.Range("N136").Value = -10858.7142
